$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (currently "un_franzosa_ControlvsCD_Fp")
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data
$ws.Range("A8").Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0.1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.4
$ws.Range("F8").Value = 0.9
$ws.Range("G8").Value = 0.5
$ws.Range("H8").Value = 0.6

# Insert a new row before what is now row 13 (currently "un_franzosa_ControlvsUC_Fp")
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new data
$ws.Range("A13").Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.4
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.6
$ws.Range("H13").Value = 0.6
